$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 19; $row++) {
    $cell = $ws.Cells.Item($row, 6)
    $current = $cell.Value2
    $newValue = $current -replace '^sequence/run_0647_samples/', ''
    $cell.Value = $newValue
}

$ws.Range("F20").Select()
